$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row stays the same text values
$ws.Range("A1").Value = "latitude"
$ws.Range("B1").Value = "longitude"

# Updated / new coordinate rows (rows 2-12)
$data = @(
    @(33.77757633306333, -84.399890863953885),
    @(33.777620190061235, -84.399819713374711),
    @(33.777615538562522, -84.399861284499607),
    @(33.777547759552618, -84.399821312264123),
    @(33.777569023561483, -84.399850092273681),
    @(33.77801490203202, -84.399895660622136),
    @(33.778015566528744, -84.399820512819417),
    @(33.778034836931781, -84.399835702268902),
    @(33.778028856462335, -84.399836501713622),
    @(33.777994967127611, -84.39983890004774),
    @(33.777650092546949, -84.399848493384255)
)

$row = 2
foreach ($pair in $data) {
    $ws.Cells.Item($row, 1).Value = $pair[0]
    $ws.Cells.Item($row, 2).Value = $pair[1]
    $row = $row + 1
}
